# ---- Update Volume number and reporting week date range text ----
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---- Update weekly crime statistics table (rows 14-29) ----
# Row 14
$ws.Range("F14").Value = 2

# Row 15
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -15.789473684210
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -40.740740740740

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -25
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = -13.513513513513
$ws.Range("I16").Value = 151
$ws.Range("J16").Value = 136
$ws.Range("K16").Value = 11.029411764705
$ws.Range("L16").Value = 32.456140350877
$ws.Range("M16").Value = 4.137931034482
$ws.Range("N16").Value = -71.292775665399

# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 15.384615384615
$ws.Range("F17").Value = 54
$ws.Range("H17").Value = 5.882352941176
$ws.Range("I17").Value = 268
$ws.Range("J17").Value = 237
$ws.Range("K17").Value = 13.080168776371
$ws.Range("L17").Value = 19.642857142857
$ws.Range("M17").Value = 79.865771812080
$ws.Range("N17").Value = -3.249097472924

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -34.782608695652
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = 4.761904761904
$ws.Range("L18").Value = 35.802469135802
$ws.Range("M18").Value = 0.917431192660
$ws.Range("N18").Value = -84.308131241084

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -5.882352941176
$ws.Range("F19").Value = 66
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -9.589041095890
$ws.Range("I19").Value = 281
$ws.Range("J19").Value = 284
$ws.Range("K19").Value = -1.056338028169
$ws.Range("L19").Value = 75.625
$ws.Range("M19").Value = 189.690721649485
$ws.Range("N19").Value = 36.407766990291

# Row 20
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 250
$ws.Range("I20").Value = 250
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = 31.578947368421
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -57.698815566835

# Row 21
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 238
$ws.Range("G21").Value = 214
$ws.Range("H21").Value = 11.214953271028
$ws.Range("I21").Value = 1081
$ws.Range("J21").Value = 988
$ws.Range("K21").Value = 9.412955465587
$ws.Range("L21").Value = 37.008871989860
$ws.Range("M21").Value = 66.563944530046
$ws.Range("N21").Value = -53.823152498932

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -62.5
$ws.Range("M22").Value = -60

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 37
$ws.Range("J23").Value = 33
$ws.Range("K23").Value = 12.121212121212
$ws.Range("L23").Value = 5.714285714285
$ws.Range("M23").Value = 117.647058823529

# Row 24
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -36.363636363636
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 177
$ws.Range("H24").Value = -36.723163841807
$ws.Range("I24").Value = 480
$ws.Range("J24").Value = 565
$ws.Range("K24").Value = -15.044247787610
$ws.Range("L24").Value = 29.380053908355
$ws.Range("M24").Value = 115.2466367713

# Row 25
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 20.833333333333
$ws.Range("F25").Value = 80
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 329
$ws.Range("J25").Value = 319
$ws.Range("K25").Value = 3.134796238244
$ws.Range("L25").Value = 20.955882352941
$ws.Range("M25").Value = -7.323943661971

# Row 26
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = -37.5
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 31
$ws.Range("K26").Value = -25.806451612903
$ws.Range("L26").Value = 4.545454545454

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 11.111111111111

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 16
$ws.Range("K28").Value = -27.272727272727
$ws.Range("L28").Value = -23.809523809523
$ws.Range("M28").Value = -20
$ws.Range("N28").Value = -65.217391304347

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 13
$ws.Range("K29").Value = -38.095238095238
$ws.Range("L29").Value = -27.777777777777
$ws.Range("M29").Value = -23.529411764705
$ws.Range("N29").Value = -69.767441860465

# ---- Fix number formats/styles for cells whose type changed between text and numeric ----
# Use PasteSpecial(Formats) from a stable donor cell of the same style to mirror
# the original workbook's conditional "0"/"***.*" vs numeric formatting.

# Donor cells with stable styles (never change type across this edit):
# D14 -> style 14 (text placeholder "0"/General, right-aligned)
# F16 -> style 15 (#,##0 numeric)
# E16 -> style 16 (signed numeric with parentheses)

$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$excel.CutCopyMode = 0